# FormatAttributeMapper.xlsx - "Optimization Updates - Attribute Mapper"
#
# A new lookup row ("other" -> "other") is inserted into the sorted
# original_value/real_value_english table at row 25 (alphabetically
# between "openapi-json" and "pdf"). This pushes the previous rows
# 25-38 down to 26-39, growing the used range from A1:B38 to A1:B39.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 25 ("pdf"/"PDF"); this shifts
# rows 25:38 down to 26:39 and preserves all existing cell content.
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row.
$ws.Range("A25").Value = "other"
$ws.Range("B25").Value = "other"

# Restore/refresh the view's active cell at the (now last) row 39 so the
# selection follows the same logical row as before the insert, matching
# a single-cell selection rather than the whole-row selection.
$ws.Range("A39").Select()

# Best-effort: scroll the window so row 28 is the top visible row
# (matches the authored sheetView topLeftCell="A28"). Some hosts may not
# persist window scroll position outside of freeze/split panes.
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
